$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.095.77'
$ws.Range('E2').Value = '  -2.08%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.991.50'
$ws.Range('E3').Value = '  -5.50%  '

$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.59'
$ws.Range('E5').Value = '  -3.10%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.21'
$ws.Range('E6').Value = '  -5.23%  '

$ws.Range('E7').Value = '  +0.17%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.986.50'
$ws.Range('E8').Value = '  -5.58%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.494'
$ws.Range('E9').Value = '  -2.49%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.134'
$ws.Range('E10').Value = '  -4.54%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.11'
$ws.Range('E11').Value = '  -2.95%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.430'
$ws.Range('E12').Value = '  -5.22%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000224'
$ws.Range('E13').Value = '  -4.02%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.98'
$ws.Range('E14').Value = '  -0.94%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.118'
$ws.Range('E15').Value = '  -0.39%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.490.97'
$ws.Range('E16').Value = '  -5.39%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.188.88'
$ws.Range('E17').Value = '  -1.92%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.004.52'
$ws.Range('E18').Value = '  -5.05%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.22'
$ws.Range('E19').Value = '  -4.62%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '434.76'
$ws.Range('E20').Value = '  -3.68%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.17'
$ws.Range('E21').Value = '  -5.58%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.662'
$ws.Range('E22').Value = '  -5.40%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.14'
$ws.Range('E23').Value = '  -5.94%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.07'
$ws.Range('E24').Value = '  -5.18%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.56'
$ws.Range('E25').Value = '  -5.28%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.19%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.03%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.50'
$ws.Range('E28').Value = '  -6.51%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.27'
$ws.Range('E29').Value = '  -5.64%  '

$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.94'
$ws.Range('E30').Value = '  -8.65%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.22'
$ws.Range('E31').Value = '  -8.65%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.32'
$ws.Range('E32').Value = '  -6.68%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0943'
$ws.Range('E33').Value = '  -8.39%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.28'
$ws.Range('E34').Value = '  -4.14%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.946'
$ws.Range('E35').Value = '  -8.32%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.56'
$ws.Range('E36').Value = '  -6.19%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '50.08'
$ws.Range('E37').Value = '  -2.18%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0678'
$ws.Range('E38').Value = '  -2.72%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0365'
$ws.Range('E39').Value = '  -4.65%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.72'
$ws.Range('E40').Value = '  -3.51%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.108'
$ws.Range('E41').Value = '  -2.26%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '371.81'
$ws.Range('E42').Value = '  -7.00%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.46'
$ws.Range('E43').Value = '  -10.54%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.645.30'
$ws.Range('E44').Value = '  -5.46%  '

$ws.Range('E45').Value = '  -0.05%  '

$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.236'
$ws.Range('E46').Value = '  -5.28%  '

$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.08'
$ws.Range('E47').Value = '  -3.08%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '33.49'
$ws.Range('E48').Value = '  -5.80%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.96'
$ws.Range('E49').Value = '  -7.38%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.106'
$ws.Range('E50').Value = '  -3.51%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.50'
$ws.Range('E51').Value = '  -6.95%  '
